$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(3.230985683306322, 1.667794583268128, 3.900430680208489, 0.496779210170732, 0, 9.295990156953671)
    3  = @(0.04763786555579896, 0.04240448674262143, 0.1575252929769615, 0.496779210170732, 0, 0.7443468554461139)
    4  = @(0.127881588408715, 0.04240448674262143, 0.8054896365839992, 0.496779210170732, 0, 1.472554921906068)
    5  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 6.201049113329182)
    6  = @(0.6753301551942219, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 0, 2.997429241610044)
    7  = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 0, 5.553084769722144)
    8  = @(0.01514828764759746, 1.667794583268128, 3.900430680208489, 8.660232485948974, 1, 14.24360603707319)
    9  = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 1, 6.201049113329182)
    10 = @(0.3048080303191223, 0.3127903958511391, 0.1575252929769615, 0.496779210170732, 1, 1.271902929317955)
    11 = @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732, 0, 6.201049113329182)
    12 = @(0.3048080303191223, 1.667794583268128, 3.900430680208489, 0.496779210170732, 1, 6.369812503966472)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 2  # B=2 ... G=7
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
